# The 3T segmentation regions (region_id 1 and 2) were swapped relative to
# the 1.5T segmentation. Fix by swapping the "temperature" (D) and
# "uncertainty" (E) values between each pair of rows that share the same
# "run"/"time" (i.e. row r and row r+1, for r = 2,4,6,...,16).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

for ($row = 2; $row -le 16; $row += 2) {
    $row2 = $row + 1

    $d1 = $ws.Cells.Item($row, 4).Value2
    $e1 = $ws.Cells.Item($row, 5).Value2
    $d2 = $ws.Cells.Item($row2, 4).Value2
    $e2 = $ws.Cells.Item($row2, 5).Value2

    $ws.Cells.Item($row, 4).Value2 = $d2
    $ws.Cells.Item($row, 5).Value2 = $e2
    $ws.Cells.Item($row2, 4).Value2 = $d1
    $ws.Cells.Item($row2, 5).Value2 = $e1
}
